$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column G and H: same new value applies to all rows 2-6
$ws.Range("G2:G6").Value = 0.56251
$ws.Range("H2:H6").Value = 1.68753

# Row 2 updates
$ws.Range("M2").Value = 20.92735966666666
$ws.Range("N2").Value = 62.782079
$ws.Range("O2").Value = 0.717329969634113
$ws.Range("P2").Value = 0.717329969634113
$ws.Range("Q2").Value = 11.77184908609666
$ws.Range("R2").Value = 105.94664177487
$ws.Range("S2").Value = 0.717329969634113
$ws.Range("T2").Value = 0.717329969634113

# Row 3 updates
$ws.Range("O3").Value = 0.009710610016949358
$ws.Range("P3").Value = 0.009710610016949358
$ws.Range("Q3").Value = 0.15935739547
$ws.Range("R3").Value = 1.43421655923
$ws.Range("S3").Value = 0.009710610016949358
$ws.Range("T3").Value = 0.009710610016949358

# Row 4 updates
$ws.Range("M4").Value = 2.195310666666666
$ws.Range("N4").Value = 6.585932
$ws.Range("O4").Value = 0.07524896398496668
$ws.Range("P4").Value = 0.07524896398496668
$ws.Range("Q4").Value = 1.234884203106666
$ws.Range("R4").Value = 11.11395782796
$ws.Range("S4").Value = 0.07524896398496668
$ws.Range("T4").Value = 0.07524896398496668

# Row 5 updates
$ws.Range("M5").Value = 3.329509666666667
$ws.Range("N5").Value = 9.988529
$ws.Range("O5").Value = 0.1141260582380437
$ws.Range("P5").Value = 0.1141260582380437
$ws.Range("Q5").Value = 1.872882482596667
$ws.Range("R5").Value = 16.85594234337
$ws.Range("S5").Value = 0.1141260582380437
$ws.Range("T5").Value = 0.1141260582380437

# Row 6 updates
$ws.Range("M6").Value = 2.438488333333333
$ws.Range("N6").Value = 7.315465
$ws.Range("O6").Value = 0.08358439812592726
$ws.Range("P6").Value = 0.08358439812592725
$ws.Range("Q6").Value = 1.371674072383333
$ws.Range("R6").Value = 12.34506665145
$ws.Range("S6").Value = 0.08358439812592726
$ws.Range("T6").Value = 0.08358439812592725
